# Commit: "complete code for project workbook" - rename the "Project" sheet
# to "Projects" and update its saved selection/active cell.

$wb = $excel.ActiveWorkbook

# Locate the sheet to rename. Prefer matching by its known original name,
# but fall back to the first sheet (its position in the workbook) if that
# name can't be found, e.g. if it was already renamed.
$wsProjects = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "Project") {
        $wsProjects = $sheet
        break
    }
}
if ($wsProjects -eq $null) {
    $wsProjects = $wb.Worksheets.Item(1)
}

$wsProjects.Name = "Projects"

# Make it the active sheet and move the selection/active cell from D8 to B9,
# matching the saved sheet view in the updated workbook.
$wsProjects.Activate()
$wsProjects.Range("B9").Select()
